$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Append new row 32 with the new inbound mail log entry
$ws.Range("A32").Value = "Afmelding nieuwsbrief"
$ws.Range("B32").Value = "mailmind.test@zohomail.eu"
$ws.Range("C32").Value = "Graag afmelden voor de nieuwsbrief. Dank u."
$ws.Range("D32").Value = "Afmelding / Nieuwsbrief"
$ws.Range("F32").Value = "2025-06-19 21:52:10"
$ws.Range("G32").Value = "Nee"

# Extend the conditional formatting ranges to cover the new row
$fcsCategorie = $ws.Range("D2:D31").FormatConditions
$fcsCategorie.Item(1).ModifyAppliesToRange($ws.Range("D2:D32"))

$fcsBeantwoord = $ws.Range("G2:G31").FormatConditions
$fcsBeantwoord.Item(1).ModifyAppliesToRange($ws.Range("G2:G32"))

# Update the Dashboard summary count for "Afmelding / Nieuwsbrief"
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B4").Value = 5
